$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.697.76"
$ws.Range("E2").Value = "  -2.78%  "
$ws.Range("D3").Value = "2.896.64"
$ws.Range("E3").Value = "  -4.06%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.45"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.35"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.503"
$ws.Range("E8").Value = "  -3.01%  "
$ws.Range("D9").Value = "2.893.16"
$ws.Range("E9").Value = "  -4.12%  "
$ws.Range("E10").Value = "  +4.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.143"
$ws.Range("E11").Value = "  -4.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("E12").Value = "  -2.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000223"
$ws.Range("E13").Value = "  -4.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.82"
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "3.377.50"
$ws.Range("E16").Value = "  -4.07%  "
$ws.Range("D17").Value = "60.683.16"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.78"
$ws.Range("E18").Value = "  -3.47%  "
$ws.Range("D19").Value = "2.900.66"
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "424.76"
$ws.Range("E20").Value = "  -5.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.56"
$ws.Range("E21").Value = "  -4.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.667"
$ws.Range("E22").Value = "  -3.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.05"
$ws.Range("E23").Value = "  -5.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.18"
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.93"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("E26").Value = "  -2.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.83"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.18"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.18"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.61"
$ws.Range("E32").Value = "  -3.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.35"
$ws.Range("E33").Value = "  -4.45%  "
$ws.Range("E34").Value = "  -3.75%  "
$ws.Range("D35").Value = "0.0₃0830"
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.63"
$ws.Range("E37").Value = "  -3.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.22"
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.02"
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.67"
$ws.Range("E42").Value = "  -4.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.288"
$ws.Range("E43").Value = "  +1.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.74"
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0345"
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "371.43"
$ws.Range("E46").Value = "  -5.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "133.72"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "2.652.80"
$ws.Range("E48").Value = "  -3.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.70"
$ws.Range("E50").Value = "  +4.03%  "
$ws.Range("E51").Value = "  -1.73%  "
